{"js": "// This document contains a table of two-digit-by-two-digit multiplication\n// problems (e.g. \"23\u00d712=276\"). The edit replaces each populated cell's\n// text with a newly generated multiplication problem/answer, keeping the\n// surrounding formatting (font, size, paragraph alignment) untouched.\nconst replacements = [\n  [\"23\u00d712=276\", \"74\u00d761=4514\"],\n  [\"89\u00d794=8366\", \"86\u00d794=8084\"],\n  [\"48\u00d723=1104\", \"35\u00d792=3220\"],\n  [\"46\u00d796=4416\", \"98\u00d758=5684\"],\n  [\"82\u00d789=7298\", \"65\u00d754=3510\"],\n  [\"67\u00d717=1139\", \"58\u00d781=4698\"],\n  [\"98\u00d727=2646\", \"53\u00d796=5088\"],\n  [\"66\u00d733=2178\", \"88\u00d756=4928\"],\n  [\"21\u00d730=630\", \"39\u00d742=1638\"],\n  [\"22\u00d727=594\", \"86\u00d779=6794\"],\n  [\"80\u00d781=6480\", \"13\u00d795=1235\"],\n  [\"46\u00d769=3174\", \"19\u00d794=1786\"],\n  [\"45\u00d784=3780\", \"82\u00d740=3280\"],\n  [\"91\u00d728=2548\", \"48\u00d753=2544\"],\n  [\"36\u00d711=396\", \"30\u00d725=750\"],\n  [\"96\u00d760=5760\", \"34\u00d794=3196\"],\n  [\"37\u00d778=2886\", \"84\u00d778=6552\"],\n  [\"32\u00d774=2368\", \"75\u00d760=4500\"],\n  [\"24\u00d767=1608\", \"33\u00d768=2244\"],\n  [\"14\u00d737=518\", \"75\u00d763=4725\"],\n  [\"98\u00d782=8036\", \"84\u00d769=5796\"],\n  [\"94\u00d734=3196\", \"55\u00d762=3410\"],\n  [\"60\u00d781=4860\", \"53\u00d723=1219\"],\n  [\"92\u00d754=4968\", \"16\u00d714=224\"],\n  [\"36\u00d718=648\", \"25\u00d760=1500\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# This document contains a table of two-digit-by-two-digit multiplication\n# problems (e.g. \"23\u00d712=276\"). The edit replaces each populated cell's\n# text with a newly generated multiplication problem/answer, keeping the\n# surrounding formatting (font, size, paragraph alignment) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"23\u00d712=276\", \"74\u00d761=4514\"),\n    @(\"89\u00d794=8366\", \"86\u00d794=8084\"),\n    @(\"48\u00d723=1104\", \"35\u00d792=3220\"),\n    @(\"46\u00d796=4416\", \"98\u00d758=5684\"),\n    @(\"82\u00d789=7298\", \"65\u00d754=3510\"),\n    @(\"67\u00d717=1139\", \"58\u00d781=4698\"),\n    @(\"98\u00d727=2646\", \"53\u00d796=5088\"),\n    @(\"66\u00d733=2178\", \"88\u00d756=4928\"),\n    @(\"21\u00d730=630\", \"39\u00d742=1638\"),\n    @(\"22\u00d727=594\", \"86\u00d779=6794\"),\n    @(\"80\u00d781=6480\", \"13\u00d795=1235\"),\n    @(\"46\u00d769=3174\", \"19\u00d794=1786\"),\n    @(\"45\u00d784=3780\", \"82\u00d740=3280\"),\n    @(\"91\u00d728=2548\", \"48\u00d753=2544\"),\n    @(\"36\u00d711=396\", \"30\u00d725=750\"),\n    @(\"96\u00d760=5760\", \"34\u00d794=3196\"),\n    @(\"37\u00d778=2886\", \"84\u00d778=6552\"),\n    @(\"32\u00d774=2368\", \"75\u00d760=4500\"),\n    @(\"24\u00d767=1608\", \"33\u00d768=2244\"),\n    @(\"14\u00d737=518\", \"75\u00d763=4725\"),\n    @(\"98\u00d782=8036\", \"84\u00d769=5796\"),\n    @(\"94\u00d734=3196\", \"55\u00d762=3410\"),\n    @(\"60\u00d781=4860\", \"53\u00d723=1219\"),\n    @(\"92\u00d754=4968\", \"16\u00d714=224\"),\n    @(\"36\u00d718=648\", \"25\u00d760=1500\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
